$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$s.Shapes.Item(5).TextFrame.TextRange.Text = "Version 2"
